# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.929.49'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '2.236.15'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.23'
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.73'
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.571'
$ws.Range("E7").Value = '  -2.28%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -4.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.46'
$ws.Range("E10").Value = '  -2.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0823'
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.39'
$ws.Range("E12").Value = '  -4.35%  '
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").Value = '2.578.96'
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.841'
$ws.Range("E15").Value = '  -2.87%  '
$ws.Range("D16").Value = '2.234.61'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("D18").Value = '43.817.21'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.90'
$ws.Range("E19").Value = '  -9.32%  '
$ws.Range("D20").Value = '0.0₃0966'
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.37'
$ws.Range("E21").Value = '  -4.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.96'
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("E23").Value = '  -3.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.99'
$ws.Range("E24").Value = '  -1.12%  '
$ws.Range("E25").Value = '  -6.31%  '
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.98'
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("E30").Value = '  -5.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.36'
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.98'
$ws.Range("E32").Value = '  -1.59%  '
$ws.Range("E33").Value = '  -4.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.68'
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.16'
$ws.Range("E35").Value = '  -2.31%  '
$ws.Range("E36").Value = '  +5.37%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("E38").Value = '  -2.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.94'
$ws.Range("E39").Value = '  +7.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.63'
$ws.Range("E40").Value = '  -2.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.03'
$ws.Range("E41").Value = '  -8.28%  '
$ws.Range("E42").Value = '  -4.11%  '
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").Value = '1.728.89'
$ws.Range("E44").Value = '  -4.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.195'
$ws.Range("E45").Value = '  -4.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '80.77'
$ws.Range("E46").Value = '  -3.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '73.52'
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("E48").Value = '  -3.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.64'
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("E51").Value = '  -3.63%  '
